# Added getMin in ForceOutcome2
# The underlying source data rows (symbol, reel1..reel5) got re-sorted
# according to a newly introduced "min" computation, so the values in
# rows 2-25 (columns A:F) of Sheet1 need to be reshuffled into the new
# row order while row 26 (the column totals) stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, 1202, 2, 10, 10, 10, 10),
    @(3, 501, 9, 52, 30, 75, 45),
    @(4, 201, 9, 30, 15, 45, 30),
    @(5, 801, 3, 67, 65, 52, 45),
    @(6, 701, 3, 90, 45, 97, 15),
    @(7, 301, 6, 45, 30, 60, 45),
    @(8, 1001, 18, 30, 75, 60, 72),
    @(9, 401, 9, 48, 67, 75, 45),
    @(10, 601, 9, 60, 67, 60, 42),
    @(11, 1201, 2, 10, 10, 10, 10),
    @(12, 1203, 3, 15, 15, 15, 15),
    @(13, 101, 9, 30, 15, 60, 15),
    @(14, 901, 16, 15, 45, 60, 60),
    @(15, 902, 1, 0, 0, 0, 0),
    @(16, 1, 0, 2, 2, 2, 2),
    @(17, 502, 0, 4, 0, 0, 0),
    @(18, 802, 0, 4, 5, 4, 0),
    @(19, 1101, 0, 15, 30, 30, 0),
    @(20, 2, 0, 2, 2, 2, 2),
    @(21, 3, 0, 3, 3, 3, 3),
    @(22, 602, 0, 0, 4, 0, 9),
    @(23, 402, 0, 0, 4, 0, 0),
    @(24, 702, 0, 0, 0, 4, 0),
    @(25, 1002, 0, 0, 0, 0, 9)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
    $ws.Cells.Item($r, 5).Value = $entry[5]
    $ws.Cells.Item($r, 6).Value = $entry[6]
}

Write-Host "Reordered rows 2-25 for ForceOutcome2 getMin sort"
